$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '35.431.42'
$ws.Cells.Item(2, 5).Value = '  -0.10%  '
$ws.Cells.Item(3, 4).Value = '1.909.78'
$ws.Cells.Item(3, 5).Value = '  +1.60%  '
$c = $ws.Cells.Item(4, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = $origStyle
$ws.Cells.Item(4, 5).Value = '  -0.46%  '
$c = $ws.Cells.Item(5, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '246.76'
$c.Style = $origStyle
$ws.Cells.Item(5, 5).Value = '  +2.56%  '
$c = $ws.Cells.Item(6, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.663'
$c.Style = $origStyle
$ws.Cells.Item(6, 5).Value = '  +6.13%  '
$c = $ws.Cells.Item(7, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.01'
$c.Style = $origStyle
$ws.Cells.Item(7, 5).Value = '  -0.51%  '
$c = $ws.Cells.Item(8, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '41.61'
$c.Style = $origStyle
$ws.Cells.Item(8, 5).Value = '  -3.13%  '
$ws.Cells.Item(9, 5).Value = '  +4.25%  '
$c = $ws.Cells.Item(10, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '53.17'
$c.Style = $origStyle
$ws.Cells.Item(10, 5).Value = '  +13.45%  '
$ws.Cells.Item(11, 5).Value = '  +2.96%  '
$ws.Cells.Item(12, 5).Value = '  -0.01%  '
$ws.Cells.Item(13, 4).Value = '2.185.00'
$ws.Cells.Item(13, 5).Value = '  +1.67%  '
$c = $ws.Cells.Item(14, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '12.27'
$c.Style = $origStyle
$ws.Cells.Item(14, 5).Value = '  +5.10%  '
$c = $ws.Cells.Item(15, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.698'
$c.Style = $origStyle
$ws.Cells.Item(15, 5).Value = '  +1.64%  '
$ws.Cells.Item(16, 4).Value = '1.915.96'
$ws.Cells.Item(16, 5).Value = '  +3.40%  '
$ws.Cells.Item(17, 5).Value = '  +1.48%  '
$ws.Cells.Item(18, 4).Value = '35.401.09'
$ws.Cells.Item(18, 5).Value = '  -0.16%  '
$c = $ws.Cells.Item(19, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '72.09'
$c.Style = $origStyle
$ws.Cells.Item(19, 5).Value = '  +2.23%  '
$ws.Cells.Item(20, 5).Value = '  +2.18%  '
$c = $ws.Cells.Item(21, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '241.34'
$c.Style = $origStyle
$ws.Cells.Item(21, 5).Value = '  -0.58%  '
$c = $ws.Cells.Item(22, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '12.47'
$c.Style = $origStyle
$ws.Cells.Item(22, 5).Value = '  +0.60%  '
$ws.Cells.Item(23, 5).Value = '  +1.21%  '
$ws.Cells.Item(24, 5).Value = '  -0.37%  '
$c = $ws.Cells.Item(25, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.41'
$c.Style = $origStyle
$ws.Cells.Item(25, 5).Value = '  +25.09%  '
$c = $ws.Cells.Item(26, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.28'
$c.Style = $origStyle
$ws.Cells.Item(26, 5).Value = '  +1.10%  '
$ws.Cells.Item(27, 5).Value = '  +0.08%  '
$c = $ws.Cells.Item(28, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '8.44'
$c.Style = $origStyle
$ws.Cells.Item(28, 5).Value = '  +2.03%  '
$c = $ws.Cells.Item(29, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '18.39'
$c.Style = $origStyle
$ws.Cells.Item(29, 5).Value = '  +3.01%  '
$ws.Cells.Item(30, 5).Value = '  +1.25%  '
$ws.Cells.Item(31, 4).Value = '4.158.48'
$ws.Cells.Item(31, 5).Value = '  +21.80%  '
$c = $ws.Cells.Item(32, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.15'
$c.Style = $origStyle
$ws.Cells.Item(32, 5).Value = '  +2.37%  '
$ws.Cells.Item(33, 2).Value = 'Hedera'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$c = $ws.Cells.Item(33, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0565'
$c.Style = $origStyle
$ws.Cells.Item(33, 5).Value = '  -0.06%  '
$ws.Cells.Item(34, 2).Value = 'ImmutableX'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$c = $ws.Cells.Item(34, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.946'
$c.Style = $origStyle
$ws.Cells.Item(34, 5).Value = '  +14.03%  '
$ws.Cells.Item(35, 5).Value = '  -0.38%  '
$c = $ws.Cells.Item(36, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '4.10'
$c.Style = $origStyle
$ws.Cells.Item(36, 5).Value = '  +0.52%  '
$c = $ws.Cells.Item(37, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.74'
$c.Style = $origStyle
$ws.Cells.Item(37, 5).Value = '  -3.10%  '
$c = $ws.Cells.Item(38, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = $origStyle
$ws.Cells.Item(38, 5).Value = '  -0.79%  '
$ws.Cells.Item(39, 5).Value = '  +1.98%  '
$c = $ws.Cells.Item(40, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '1.11'
$c.Style = $origStyle
$ws.Cells.Item(40, 5).Value = '  -1.40%  '
$ws.Cells.Item(41, 5).Value = '  +2.43%  '
$ws.Cells.Item(42, 2).Value = 'Kaspa'
$ws.Cells.Item(42, 3).Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$c = $ws.Cells.Item(42, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '0.0647'
$c.Style = $origStyle
$ws.Cells.Item(42, 5).Value = '  +7.69%  '
$ws.Cells.Item(43, 2).Value = 'InjectiveProtocol'
$ws.Cells.Item(43, 3).Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$c = $ws.Cells.Item(43, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '16.26'
$c.Style = $origStyle
$ws.Cells.Item(43, 5).Value = '  +6.06%  '
$c = $ws.Cells.Item(44, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '89.83'
$c.Style = $origStyle
$ws.Cells.Item(44, 5).Value = '  -1.80%  '
$ws.Cells.Item(45, 4).Value = '1.336.44'
$ws.Cells.Item(45, 5).Value = '  -1.42%  '
$c = $ws.Cells.Item(46, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '49.36'
$c.Style = $origStyle
$ws.Cells.Item(46, 5).Value = '  +40.62%  '
$c = $ws.Cells.Item(47, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.40'
$c.Style = $origStyle
$ws.Cells.Item(47, 5).Value = '  +1.03%  '
$c = $ws.Cells.Item(48, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.78'
$c.Style = $origStyle
$ws.Cells.Item(48, 5).Value = '  +2.67%  '
$c = $ws.Cells.Item(49, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '2.39'
$c.Style = $origStyle
$ws.Cells.Item(49, 5).Value = '  -0.75%  '
$c = $ws.Cells.Item(50, 4)
$origStyle = $c.Style
$c.NumberFormat = "@"
$c.Value = '6.50'
$c.Style = $origStyle
$ws.Cells.Item(50, 5).Value = '  -2.45%  '
$ws.Cells.Item(51, 4).Value = '2.092.82'
$ws.Cells.Item(51, 5).Value = '  +1.35%  '
